# Update "countries & provincias Spain" dataset (refreshed COVID-19 snapshot).
# The underlying source data was re-sorted by "Casos totales" (col B, desc);
# a handful of countries changed rank, which re-labels some rows, and the
# statistics themselves were refreshed for this batch of rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- "Datos actualizados ..." timestamp in A1 ---
$ws.Range("A1").Value = "Datos actualizados a 25 de Julio de 2020 a las 21:50"

# --- Country column (A) re-labels from the re-sort ---
$ws.Range("A72").Value = "Costa Rica"
$ws.Range("A73").Value = "Venezuela"
$ws.Range("A74").Value = "El Salvador"
$ws.Range("A75").Value = "Corea del Sur"
$ws.Range("A76").Value = "Australia"
$ws.Range("A109").Value = "Somalia"
$ws.Range("A110").Value = "Maldivas"
$ws.Range("A148").Value = "Angola"
$ws.Range("A149").Value = "Republica del Chad"
$ws.Range("A150").Value = "Principado de Andorra"

# --- Refreshed statistics (B:Casos totales, C:Nuevos casos, D:Casos activos, E:Recuperados, G:Casos criticos, H:Muertes) ---
$ws.Range("B4").Value = 4290235
$ws.Range("C4").Value = 41908
$ws.Range("D4").Value = 2040513
$ws.Range("E4").Value = 2100675
$ws.Range("G4").Value = 557
$ws.Range("H4").Value = 149047

$ws.Range("B21").Value = 206276
$ws.Range("C21").Value = 316
$ws.Range("E21").Value = 6674
$ws.Range("G21").Value = 1
$ws.Range("H21").Value = 9202

$ws.Range("B65").Value = 19952
$ws.Range("C65").Value = 592
$ws.Range("D65").Value = 10831
$ws.Range("E65").Value = 9010
$ws.Range("G65").Value = 5
$ws.Range("H65").Value = 111

$ws.Range("B72").Value = 14600
$ws.Range("C72").Value = 931
$ws.Range("D72").Value = 3640
$ws.Range("E72").Value = 10862
$ws.Range("G72").Value = 11
$ws.Range("H72").Value = 98

$ws.Range("B73").Value = 14263
$ws.Range("C73").Value = 0
$ws.Range("D73").Value = 8127
$ws.Range("E73").Value = 6002
$ws.Range("G73").Value = 0
$ws.Range("H73").Value = 134

$ws.Range("B74").Value = 14221
$ws.Range("C74").Value = 429
$ws.Range("D74").Value = 7549
$ws.Range("E74").Value = 6282
$ws.Range("G74").Value = 11
$ws.Range("H74").Value = 390

$ws.Range("B75").Value = 14092
$ws.Range("C75").Value = 113
$ws.Range("D75").Value = 12866
$ws.Range("E75").Value = 928
$ws.Range("G75").Value = 0
$ws.Range("H75").Value = 298

$ws.Range("B76").Value = 13948
$ws.Range("C76").Value = 353
$ws.Range("D76").Value = 9017
$ws.Range("E76").Value = 4786
$ws.Range("G76").Value = 6
$ws.Range("H76").Value = 145

$ws.Range("B79").Value = 11385
$ws.Range("C79").Value = 83
$ws.Range("D79").Value = 5890
$ws.Range("E79").Value = 4778
$ws.Range("G79").Value = 2
$ws.Range("H79").Value = 717

$ws.Range("B109").Value = 3178
$ws.Range("C109").Value = 7
$ws.Range("D109").Value = 1521
$ws.Range("E109").Value = 1564
$ws.Range("H109").Value = 93

$ws.Range("B110").Value = 3175
$ws.Range("D110").Value = 2498
$ws.Range("E110").Value = 662
$ws.Range("H110").Value = 15

$ws.Range("D145").Value = 982
$ws.Range("E145").Value = 120

$ws.Range("B148").Value = 916
$ws.Range("C148").Value = 36
$ws.Range("D148").Value = 242
$ws.Range("E148").Value = 635
$ws.Range("G148").Value = 4
$ws.Range("H148").Value = 39

$ws.Range("B149").Value = 915
$ws.Range("D149").Value = 805
$ws.Range("E149").Value = 35
$ws.Range("H149").Value = 75

$ws.Range("B150").Value = 897
$ws.Range("D150").Value = 803
$ws.Range("E150").Value = 42
$ws.Range("H150").Value = 52
